$wb = $excel.ActiveWorkbook

# Rename Sheet1 -> ValidLogin
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ValidLogin"

# Fill in the login test-data table.
# Write order matters for shared-string pool ordering (first-seen wins):
# Username, Password, manager, admin, HomePageTitle, actiTIME - Enter Time-Track
$ws1.Range("A1").Value = "Username"
$ws1.Range("B1").Value = "Password"
$ws1.Range("B2").Value = "manager"
$ws1.Range("A2").Value = "admin"
$ws1.Range("C1").Value = "HomePageTitle"
$ws1.Range("C2").Value = "actiTIME - Enter Time-Track"

# Bold header row
$ws1.Range("A1:C1").Font.Bold = $true

# Column widths (best-fit-like), closest achievable to the recorded widths
$ws1.Columns.Item(1).ColumnWidth = 9.3
$ws1.Columns.Item(2).ColumnWidth = 8.6
$ws1.Columns.Item(3).ColumnWidth = 24.5

# Force a pageSetup element (portrait) to be emitted
$ws1.PageSetup.Orientation = 1

# Sheet2: remember a selection on A2, then return focus to Sheet1
$ws2 = $wb.Worksheets.Item(2)
$null = $ws2.Range("A2").Select()
$null = $ws1.Select()
